$wb = $excel.ActiveWorkbook

# Rename sheets to unify DataNode / DataTable naming
$wsNode = $wb.Worksheets.Item(1)
$wsNode.Name = "DataNode"

$wsTable = $wb.Worksheets.Item(2)
$wsTable.Name = "DataTable"

# Adjust header row heights on the DataNode sheet
$wsNode.Rows.Item(1).RowHeight = 27
$wsNode.Rows.Item(8).RowHeight = 54

# Move selection / activation to the DataTable sheet
$wsTable.Activate()
$wsTable.Range("H32").Select()
